$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4 (Kalenderwoche), Row 5 (Theke), Row 6 (Zweigstelle) get reshuffled
# across columns C:N.

$ws.Range("C4").Value = "23. KW"
$ws.Range("D4").Value = "24. KW"
$ws.Range("E4").Value = "23. KW"
$ws.Range("F4").Value = "23. KW"
$ws.Range("G4").Value = "23. KW"
$ws.Range("H4").Value = "23. KW"
$ws.Range("I4").Value = "23. KW"
$ws.Range("J4").Value = "24. KW"
$ws.Range("K4").Value = "24. KW"
$ws.Range("L4").Value = "24. KW"
$ws.Range("M4").Value = "24. KW"
$ws.Range("N4").Value = "24. KW"

$ws.Range("C5").Value = "Auskunftstheke"
$ws.Range("D5").Value = "Auskunftstheke"
$ws.Range("E5").Value = "Auskunftstheke"
$ws.Range("F5").Value = "Auskunftstheke"
$ws.Range("G5").Value = "Ausleihtheke"
$ws.Range("H5").Value = "Ausleihtheke"
$ws.Range("I5").Value = "Ausleihtheke"
$ws.Range("J5").Value = "Auskunftstheke"
$ws.Range("K5").Value = "Ausleihtheke"
$ws.Range("L5").Value = "Ausleihtheke"
$ws.Range("M5").Value = "Ausleihtheke"
$ws.Range("N5").Value = "Auskunftstheke"

$ws.Range("C6").Value = "GM"
$ws.Range("D6").Value = "Deutz"
$ws.Range("E6").Value = "Deutz"
$ws.Range("F6").Value = "Südstadt"
$ws.Range("G6").Value = "Deutz"
$ws.Range("H6").Value = "GM"
$ws.Range("I6").Value = "Südstadt"
$ws.Range("J6").Value = "GM"
$ws.Range("K6").Value = "Deutz"
$ws.Range("L6").Value = "GM"
$ws.Range("M6").Value = "Südstadt"
$ws.Range("N6").Value = "Südstadt"
